# Applies the odds updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("G6").Value = 2.32
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.72
$ws.Range("J6").Value = 2.92
$ws.Range("L6").Value = 3.25
$ws.Range("V6").Value = 2.05
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 9.25
$ws.Range("Z6").Value = 24
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 11.5
$ws.Range("AE6").Value = 13
$ws.Range("AP6").Value = 19.5
$ws.Range("AR6").Value = 80
$ws.Range("AS6").Value = 250
$ws.Range("AU6").Value = 6.8
$ws.Range("AW6").Value = 4.7
$ws.Range("AZ6").Value = 60

# Row 7
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 7
$ws.Range("K7").Value = 2.47
$ws.Range("S7").Value = 1.28
$ws.Range("T7").Value = 3.46
$ws.Range("Z7").Value = 9
$ws.Range("AC7").Value = 15.5
$ws.Range("AF7").Value = 75
$ws.Range("AG7").Value = 22
$ws.Range("AK7").Value = 70
$ws.Range("AO7").Value = 6
$ws.Range("AT7").Value = 3.3
$ws.Range("AU7").Value = 7.9
$ws.Range("AY7").Value = 37

# Row 13
$ws.Range("M13").Value = 1.07
$ws.Range("O13").Value = 1.33

# Row 14
$ws.Range("M14").Value = 1.07
$ws.Range("O14").Value = 1.33

# Row 15
$ws.Range("M15").Value = 1.1
$ws.Range("O15").Value = 1.44
$ws.Range("P15").Value = 2.63

# Row 16
$ws.Range("M16").Value = 1.07
$ws.Range("O16").Value = 1.33
